# Weekly refresh of the "Fruta / hortaliza" price data.
# The data rows (2..43) get re-associated: each destination row ends up
# holding the Fecha/Volumen/Precio values that another row used to hold.
# We snapshot all source values first, then write them back out so that
# overlapping reads/writes never clobber data we still need.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destination row -> source row (values currently sitting in source row
# should end up in destination row)
$map = @{
  2  = 40
  3  = 11
  4  = 29
  5  = 42
  6  = 23
  7  = 28
  8  = 6
  9  = 20
  10 = 25
  11 = 31
  12 = 39
  13 = 16
  14 = 7
  15 = 30
  16 = 15
  17 = 35
  18 = 4
  19 = 38
  20 = 13
  21 = 19
  22 = 33
  23 = 5
  24 = 41
  25 = 17
  26 = 21
  27 = 14
  28 = 24
  29 = 2
  30 = 32
  31 = 37
  32 = 9
  33 = 27
  34 = 22
  35 = 26
  36 = 12
  37 = 10
  38 = 18
  39 = 34
  40 = 36
  41 = 8
  42 = 43
  43 = 3
}

# Columns that move together with each data row: D=Fecha, J=Volumen,
# K=Precio minimo, L=Precio maximo, M=Precio promedio ponderado,
# P=Precio $/Kg (always mirrors M in this sheet).
$cols = @(4, 10, 11, 12, 13, 16)

# Snapshot current values for every row/column we touch before writing
# anything, since a row can be both a source and a destination.
$snapshot = @{}
for ($row = 2; $row -le 43; $row++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($row, $c).Value2
    }
    $snapshot[$row] = $rowVals
}

for ($row = 2; $row -le 43; $row++) {
    $srcRow = $map[$row]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($row, $c).Value2 = $srcVals[$c]
    }
}
